$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column S, year 2022 (copy formatting from R4 which has the same style)
$ws.Range("R4").Copy()
$ws.Range("S4").PasteSpecial(-4122)
$ws.Range("S4").Value = 2022

# Update existing values in row 5
$ws.Range("P5").Value = 20.5
$ws.Range("Q5").Value = 20.5
$ws.Range("R5").Value = 17.899999999999999

# New value for S5 (copy formatting from R5 which has the same style)
$ws.Range("R5").Copy()
$ws.Range("S5").PasteSpecial(-4122)
$ws.Range("S5").Value = 13.5

$excel.CutCopyMode = $false

# Update selection to match diff
$ws.Range("S7:S8").Select()
